# Update column G ("K") values on Sheet1 for rows 2-17.
# These values represent strikeouts (Strike#) being replaced by a
# recalculated "K" stat as part of regenerating save_data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 2
    3  = 2
    4  = 1
    5  = 0
    6  = 8
    7  = 3
    8  = 5
    9  = 1
    10 = 3
    11 = 6
    12 = 2
    13 = 3
    14 = 5
    15 = 3
    16 = 1
    17 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
